$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new BOM line item for the proc board connector
$ws.Cells.Item(32, 1).Value = 1
$ws.Cells.Item(32, 2).Value = "A115899CT-ND"
$ws.Cells.Item(32, 3).Value = "PROC BOARD CONNECTOR"
$ws.Cells.Item(32, 4).Value = "CONN FEMALE 67POS 0.020 GOLD"

# Reflect where the user was scrolled/selected after making the edit
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.Zoom = 100
$ws.Range("D33").Select()
